# Update cryptos list - price and 1h volume % changes scraped on
# Sun Nov 17 10:42:22 UTC 2024 via GitHub Actions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.091.76"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.142.40"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +9.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.58"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.366"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.139.76"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.722"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.00%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.05"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.55"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.913.88"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.721.03"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.090.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000212"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +8.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.04"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.295.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.69"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("E32").Value = "  +9.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +11.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.197"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +23.21%  "
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "515.84"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.30"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.419"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.22"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0853"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +47.42%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.66"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.695"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.68%  "

# Row 49 (was ImmutableX) becomes OKB; row 50 (was OKB) becomes ImmutableX.
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.88"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.37"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.10%  "
